# Se agrega complejidad computacional de los binomios en el excel
#
# The "Binomio de Newton" summary table (columns K:Q, header row 14 with
# Msucesiva, Recursiva, RecursivaPar, Horner, ProgDinam, Mejorado, Pow)
# had an empty "Complejidad Computacional" row (row 15). This fills in
# the Big-O complexity for each algorithm, mirroring the values already
# present for the "Polinomios" table in columns B:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("K15").Value = "n"    # Msucesiva
$ws.Range("L15").Value = "n2"   # Recursiva
$ws.Range("M15").Value = "n2"   # RecursivaPar
$ws.Range("N15").Value = "n"    # Horner
$ws.Range("O15").Value = "n"    # ProgDinam
$ws.Range("Q15").Value = "n"    # Pow
